$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Age" -> "Age " (trailing space) and make header bold
$ws.Range("A1").Value = "Age "
$ws.Range("B1").Value = "Gender"
$ws.Range("C1").Value = "Genre"
$ws.Range("A1:C1").Font.Bold = $true

# Genre label correction: ROCK -> Rock, and updated/added data rows
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "Rock"

$ws.Range("A3").Value = 22
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "Rock"

$ws.Range("A4").Value = 24
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "Rock"

$ws.Range("A5").Value = 26
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Rock"

$ws.Range("A6").Value = 27
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "JAZZ"

$ws.Range("A7").Value = 29
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "JAZZ"

$ws.Range("A8").Value = 31
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "JAZZ"

$ws.Range("A9").Value = 35
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "JAZZ"

$ws.Range("A10").Value = 36
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "CLASSICAL"

$ws.Range("A11").Value = 40
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = "CLASSICAL"

$ws.Range("A12").Value = 45
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = "CLASSICAL"

$ws.Range("A13").Value = 52
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = "CLASSICAL"

$ws.Range("A14").Value = 20
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Rock"

$ws.Range("A15").Value = 22
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Rock"

$ws.Range("A16").Value = 24
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Rock"

$ws.Range("A17").Value = 26
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Rock"

$ws.Range("A18").Value = 27
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "JAZZ"

$ws.Range("A19").Value = 29
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "JAZZ"

$ws.Range("A20").Value = 31
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "JAZZ"

$ws.Range("A21").Value = 35
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "CLASSICAL"

$ws.Range("A22").Value = 36
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "CLASSICAL"

$ws.Range("A23").Value = 40
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "CLASSICAL"

# Two new rows appended at the bottom - clone row 23's formatting (border +
# centered alignment) down onto them before writing their values.
$ws.Range("A23:C23").Copy($ws.Range("A24:C24"))
$ws.Range("A23:C23").Copy($ws.Range("A25:C25"))

$ws.Range("A24").Value = 45
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "CLASSICAL"

$ws.Range("A25").Value = 52
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "CLASSICAL"

# View changes: zoom to 122% and move the (non-data) selection to F8
$ws.Application.ActiveWindow.Zoom = 122
$ws.Range("F8").Select() | Out-Null
